$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 146.620486
$ws.Range("H2").Value = 439.861458
$ws.Range("I2").Value = 0.3983053592962091
$ws.Range("J2").Value = 0.3983053592962091
$ws.Range("M2").Value = 35.04689966666667
$ws.Range("N2").Value = 105.140699
$ws.Range("O2").Value = 0.3824629895491901
$ws.Range("P2").Value = 0.3824629895491901
$ws.Range("Q2").Value = 5138.593461919905
$ws.Range("R2").Value = 46247.34115727914
$ws.Range("S2").Value = 0.1523370584698924
$ws.Range("T2").Value = 0.1523370584698924
$ws.Range("G3").Value = 146.620486
$ws.Range("H3").Value = 439.861458
$ws.Range("I3").Value = 0.3983053592962091
$ws.Range("J3").Value = 0.3983053592962091
$ws.Range("O3").Value = 0.3264402385872224
$ws.Range("P3").Value = 0.3264402385872223
$ws.Range("Q3").Value = 4385.898038628734
$ws.Range("R3").Value = 39473.0823476586
$ws.Range("S3").Value = 0.1300228965192238
$ws.Range("T3").Value = 0.1300228965192238
$ws.Range("G4").Value = 146.620486
$ws.Range("H4").Value = 439.861458
$ws.Range("I4").Value = 0.3983053592962091
$ws.Range("J4").Value = 0.3983053592962091
$ws.Range("M4").Value = 8.911727666666666
$ws.Range("N4").Value = 26.735183
$ws.Range("O4").Value = 0.09725271102035077
$ws.Range("P4").Value = 0.09725271102035075
$ws.Range("Q4").Value = 1306.641841586312
$ws.Range("R4").Value = 11759.77657427681
$ws.Range("S4").Value = 0.0387362760054912
$ws.Range("T4").Value = 0.0387362760054912
$ws.Range("G5").Value = 146.620486
$ws.Range("H5").Value = 439.861458
$ws.Range("I5").Value = 0.3983053592962091
$ws.Range("J5").Value = 0.3983053592962091
$ws.Range("M5").Value = 17.76285166666667
$ws.Range("N5").Value = 53.288555
$ws.Range("O5").Value = 0.1938440608432367
$ws.Range("P5").Value = 0.1938440608432367
$ws.Range("Q5").Value = 2604.397944112577
$ws.Range("R5").Value = 23439.58149701319
$ws.Range("S5").Value = 0.07720912830160162
$ws.Range("T5").Value = 0.07720912830160162
$ws.Range("I6").Value = 0.534552907532962
$ws.Range("J6").Value = 0.5345529075329621
$ws.Range("M6").Value = 35.04689966666667
$ws.Range("N6").Value = 105.140699
$ws.Range("O6").Value = 0.3824629895491901
$ws.Range("P6").Value = 0.3824629895491901
$ws.Range("Q6").Value = 6896.342244936743
$ws.Range("R6").Value = 62067.08020443069
$ws.Range("S6").Value = 0.2044467030872684
$ws.Range("T6").Value = 0.2044467030872684
$ws.Range("I7").Value = 0.534552907532962
$ws.Range("J7").Value = 0.5345529075329621
$ws.Range("O7").Value = 0.3264402385872224
$ws.Range("P7").Value = 0.3264402385872223
$ws.Range("S7").Value = 0.1744995786725536
$ws.Range("T7").Value = 0.1744995786725536
$ws.Range("I8").Value = 0.534552907532962
$ws.Range("J8").Value = 0.5345529075329621
$ws.Range("M8").Value = 8.911727666666666
$ws.Range("N8").Value = 26.735183
$ws.Range("O8").Value = 0.09725271102035077
$ws.Range("P8").Value = 0.09725271102035075
$ws.Range("Q8").Value = 1753.602303414538
$ws.Range("R8").Value = 15782.42073073084
$ws.Range("S8").Value = 0.05198671944139144
$ws.Range("T8").Value = 0.05198671944139144
$ws.Range("I9").Value = 0.534552907532962
$ws.Range("J9").Value = 0.5345529075329621
$ws.Range("M9").Value = 17.76285166666667
$ws.Range("N9").Value = 53.288555
$ws.Range("O9").Value = 0.1938440608432367
$ws.Range("P9").Value = 0.1938440608432367
$ws.Range("Q9").Value = 3495.279340097739
$ws.Range("R9").Value = 31457.51406087966
$ws.Range("S9").Value = 0.1036199063317486
$ws.Range("T9").Value = 0.1036199063317486
$ws.Range("G10").Value = 24.174389
$ws.Range("H10").Value = 72.523167
$ws.Range("I10").Value = 0.0656715098899026
$ws.Range("J10").Value = 0.0656715098899026
$ws.Range("M10").Value = 35.04689966666667
$ws.Range("N10").Value = 105.140699
$ws.Range("O10").Value = 0.3824629895491901
$ws.Range("P10").Value = 0.3824629895491901
$ws.Range("Q10").Value = 847.2373857859704
$ws.Range("R10").Value = 7625.136472073734
$ws.Range("S10").Value = 0.02511692200070135
$ws.Range("T10").Value = 0.02511692200070135
$ws.Range("G11").Value = 24.174389
$ws.Range("H11").Value = 72.523167
$ws.Range("I11").Value = 0.0656715098899026
$ws.Range("J11").Value = 0.0656715098899026
$ws.Range("O11").Value = 0.3264402385872224
$ws.Range("P11").Value = 0.3264402385872223
$ws.Range("Q11").Value = 723.135001067641
$ws.Range("R11").Value = 6508.215009608769
$ws.Range("S11").Value = 0.02143782335684294
$ws.Range("T11").Value = 0.02143782335684294
$ws.Range("G12").Value = 24.174389
$ws.Range("H12").Value = 72.523167
$ws.Range("I12").Value = 0.0656715098899026
$ws.Range("J12").Value = 0.0656715098899026
$ws.Range("M12").Value = 8.911727666666666
$ws.Range("N12").Value = 26.735183
$ws.Range("O12").Value = 0.09725271102035077
$ws.Range("P12").Value = 0.09725271102035075
$ws.Range("Q12").Value = 215.4355712760623
$ws.Range("R12").Value = 1938.920141484561
$ws.Range("S12").Value = 0.006386732373592805
$ws.Range("T12").Value = 0.006386732373592804
$ws.Range("G13").Value = 24.174389
$ws.Range("H13").Value = 72.523167
$ws.Range("I13").Value = 0.0656715098899026
$ws.Range("J13").Value = 0.0656715098899026
$ws.Range("M13").Value = 17.76285166666667
$ws.Range("N13").Value = 53.288555
$ws.Range("O13").Value = 0.1938440608432367
$ws.Range("P13").Value = 0.1938440608432367
$ws.Range("Q13").Value = 429.4060859392984
$ws.Range("R13").Value = 3864.654773453685
$ws.Range("S13").Value = 0.0127300321587655
$ws.Range("T13").Value = 0.0127300321587655
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.541205
$ws.Range("H14").Value = 1.623615
$ws.Range("I14").Value = 0.001470223280926138
$ws.Range("J14").Value = 0.001470223280926138
$ws.Range("M14").Value = 35.04689966666667
$ws.Range("N14").Value = 105.140699
$ws.Range("O14").Value = 0.3824629895491901
$ws.Range("P14").Value = 0.3824629895491901
$ws.Range("Q14").Value = 18.96755733409834
$ws.Range("R14").Value = 170.708016006885
$ws.Range("S14").Value = 0.0005623059913278295
$ws.Range("T14").Value = 0.0005623059913278294
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.541205
$ws.Range("H15").Value = 1.623615
$ws.Range("I15").Value = 0.001470223280926138
$ws.Range("J15").Value = 0.001470223280926138
$ws.Range("O15").Value = 0.3264402385872224
$ws.Range("P15").Value = 0.3264402385872223
$ws.Range("Q15").Value = 16.189210749145
$ws.Range("R15").Value = 145.702896742305
$ws.Range("S15").Value = 0.0004799400386020174
$ws.Range("T15").Value = 0.0004799400386020172
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.541205
$ws.Range("H16").Value = 1.623615
$ws.Range("I16").Value = 0.001470223280926138
$ws.Range("J16").Value = 0.001470223280926138
$ws.Range("M16").Value = 8.911727666666666
$ws.Range("N16").Value = 26.735183
$ws.Range("O16").Value = 0.09725271102035077
$ws.Range("P16").Value = 0.09725271102035075
$ws.Range("Q16").Value = 4.823071571838334
$ws.Range("R16").Value = 43.407644146545
$ws.Range("S16").Value = 0.0001429831998753017
$ws.Range("T16").Value = 0.0001429831998753016
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.541205
$ws.Range("H17").Value = 1.623615
$ws.Range("I17").Value = 0.001470223280926138
$ws.Range("J17").Value = 0.001470223280926138
$ws.Range("M17").Value = 17.76285166666667
$ws.Range("N17").Value = 53.288555
$ws.Range("O17").Value = 0.1938440608432367
$ws.Range("P17").Value = 0.1938440608432367
$ws.Range("Q17").Value = 9.613344136258334
$ws.Range("R17").Value = 86.52009722632501
$ws.Range("S17").Value = 0.0002849940511209894
$ws.Range("T17").Value = 0.0002849940511209894